$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count

# Columns (1-based) that can contain the Distance/Size coded substrings:
#   B (2)  = Condition
#   D (4)  = Filename_Left
#   E (5)  = Filename_Right
#   H (8)  = Distance
#   J (10) = Size
$targetCols = @(2, 4, 5, 8, 10)

for ($r = 2; $r -le $rowCount; $r++) {
    foreach ($c in $targetCols) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -ne $null -and $val.GetType().Name -eq "String") {
            $newVal = $val.Replace("D64", "D69").Replace("D51", "D55").Replace("D80", "D86").Replace("S30", "S31")
            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
